# Auto-generated Excel COM-interop script to apply the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Test-NumberLike($s) {
    return $s -match '^[+-]?[0-9]+(\.[0-9]+)?$'
}

function Set-TextValue($range, $value) {
    if (Test-NumberLike $value) {
        $range.NumberFormat = "@"
    }
    $range.Value = $value
}

# Row 2: Bitcoin
Set-TextValue $ws.Range("D2") '28.949.94'
Set-TextValue $ws.Range("E2") '  -2.03%  '

# Row 3: Ethereum
Set-TextValue $ws.Range("D3") '1.904.50'
Set-TextValue $ws.Range("E3") '  -4.33%  '

# Row 4: TetherUSD
Set-TextValue $ws.Range("D4") '1.005'
Set-TextValue $ws.Range("E4") '  +0.10%  '

# Row 5: BNB
Set-TextValue $ws.Range("D5") '324.83'
Set-TextValue $ws.Range("E5") '  -0.15%  '

# Row 6: USDC
Set-TextValue $ws.Range("D6") '1.003'
Set-TextValue $ws.Range("E6") '  +0.22%  '

# Row 7: XRP
Set-TextValue $ws.Range("D7") '0.4596'
Set-TextValue $ws.Range("E7") '  -1.73%  '

# Row 8: Cardano
Set-TextValue $ws.Range("D8") '0.3819'
Set-TextValue $ws.Range("E8") '  -3.17%  '

# Row 9: OKB
Set-TextValue $ws.Range("E9") '  -2.15%  '

# Row 10: Dogecoin
Set-TextValue $ws.Range("D10") '0.07734'
Set-TextValue $ws.Range("E10") '  -2.67%  '

# Row 11: Polygon
Set-TextValue $ws.Range("D11") '0.9824'
Set-TextValue $ws.Range("E11") '  -2.02%  '

# Row 12: Solana
Set-TextValue $ws.Range("D12") '22.06'
Set-TextValue $ws.Range("E12") '  -4.05%  '

# Row 13: WrappedEther
Set-TextValue $ws.Range("D13") '1.925.23'
Set-TextValue $ws.Range("E13") '  -3.14%  '

# Row 14: Chainlink
Set-TextValue $ws.Range("D14") '6.992'
Set-TextValue $ws.Range("E14") '  -3.72%  '

# Row 15: Polkadot
Set-TextValue $ws.Range("E15") '  -3.37%  '

# Row 16: TRON
Set-TextValue $ws.Range("D16") '0.07041'
Set-TextValue $ws.Range("E16") '  -1.47%  '

# Row 17: BinanceUSD
Set-TextValue $ws.Range("D17") '1.005'
Set-TextValue $ws.Range("E17") '  +0.17%  '

# Row 18: Litecoin
Set-TextValue $ws.Range("D18") '84.18'
Set-TextValue $ws.Range("E18") '  -5.19%  '

# Row 19: ShibaInu
Set-TextValue $ws.Range("D19") '0.000009554'

# Row 20: Avalanche
Set-TextValue $ws.Range("D20") '16.76'
Set-TextValue $ws.Range("E20") '  -3.76%  '

# Row 21: Dai
Set-TextValue $ws.Range("E21") '  +0.27%  '

# Row 22: WrappedBTC
Set-TextValue $ws.Range("D22") '28.981.71'
Set-TextValue $ws.Range("E22") '  -2.17%  '

# Row 23: Uniswap
Set-TextValue $ws.Range("D23") '5.331'
Set-TextValue $ws.Range("E23") '  -3.93%  '

# Row 24: Cosmos
Set-TextValue $ws.Range("E24") '  -2.89%  '

# Row 25: WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D25") '2.166.81'
Set-TextValue $ws.Range("E25") '  -2.64%  '

# Row 26: Toncoin
Set-TextValue $ws.Range("D26") '2.078'
Set-TextValue $ws.Range("E26") '  -1.28%  '

# Row 27: Monero
Set-TextValue $ws.Range("D27") '156.57'
Set-TextValue $ws.Range("E27") '  -0.89%  '

# Row 28: EthereumClassic
Set-TextValue $ws.Range("D28") '19.13'
Set-TextValue $ws.Range("E28") '  -2.90%  '

# Row 29: InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") '5.594'
Set-TextValue $ws.Range("E29") '  -6.18%  '

# Row 30: BitcoinCash
Set-TextValue $ws.Range("D30") '117.72'
Set-TextValue $ws.Range("E30") '  -2.26%  '

# Row 31: LidoDAOToken
Set-TextValue $ws.Range("D31") '1.832'
Set-TextValue $ws.Range("E31") '  -7.01%  '

# Row 32: Stellar
Set-TextValue $ws.Range("D32") '0.09259'
Set-TextValue $ws.Range("E32") '  -2.05%  '

# Row 33: ImmutableX
Set-TextValue $ws.Range("D33") '0.8603'
Set-TextValue $ws.Range("E33") '  -4.57%  '

# Row 34: Filecoin
Set-TextValue $ws.Range("D34") '5.111'
Set-TextValue $ws.Range("E34") '  -2.90%  '

# Row 35: ARBITRUM
Set-TextValue $ws.Range("E35") '  -6.82%  '

# Row 36: HuobiToken
Set-TextValue $ws.Range("D36") '3.017'
Set-TextValue $ws.Range("E36") '  -5.24%  '

# Row 37: Hedera
Set-TextValue $ws.Range("D37") '0.05726'
Set-TextValue $ws.Range("E37") '  -2.03%  '

# Row 38: TrustWalletToken
Set-TextValue $ws.Range("D38") '1.145'
Set-TextValue $ws.Range("E38") '  -2.70%  '

# Row 39: Frax
Set-TextValue $ws.Range("D39") '1.004'
Set-TextValue $ws.Range("E39") '  +0.26%  '

# Row 40: VeChain
Set-TextValue $ws.Range("D40") '0.02041'
Set-TextValue $ws.Range("E40") '  -3.87%  '

# Row 41: FraxShare
Set-TextValue $ws.Range("D41") '7.499'
Set-TextValue $ws.Range("E41") '  -5.21%  '

# Row 42: TheSandbox
Set-TextValue $ws.Range("D42") '0.5532'
Set-TextValue $ws.Range("E42") '  -4.01%  '

# Row 43: Algorand
Set-TextValue $ws.Range("E43") '  -4.03%  '

# Row 44: Aptos
Set-TextValue $ws.Range("E44") '  -5.34%  '

# Row 45: MXToken
Set-TextValue $ws.Range("D45") '2.758'
Set-TextValue $ws.Range("E45") '  +2.45%  '

# Row 46: Decentraland
Set-TextValue $ws.Range("D46") '0.5217'
Set-TextValue $ws.Range("E46") '  -2.98%  '

# Row 47: EnergySwap
Set-TextValue $ws.Range("E47") '  -6.32%  '

# Row 48: RenderToken
Set-TextValue $ws.Range("D48") '2.092'
Set-TextValue $ws.Range("E48") '  -4.29%  '

# Rows 49-51: reorder PEPE / Cronos / Quant
# Row 49: Cronos -> PEPE
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D49") '0.000002634'
$ws.Range("E49").Value = '  -20.38%  '

# Row 50: Quant -> Cronos
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.06821'
$ws.Range("E50").Value = '  -1.90%  '

# Row 51: PEPE -> Quant
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D51") '111.97'
$ws.Range("E51").Value = '  -2.32%  '

Write-Host "Applied cryptos.xlsx update"